$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8580e004da770ffceef172cdbe1908c444750cf0/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5c91a3163706bbd7f1af8bdf16a792200c539785/e2e/b.md."

# --- Overview sheet: row 3 is b.md ---
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-19 20:42:54"

# --- zh-cn sheet: row 3 is b.md ---
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-19 20:42:49"
$wsZhCn.Range("P3").Value = $errorDetail
$wsZhCn.Range("P1").EntireColumn.ColumnWidth = 39.1666666666667

# --- de-de sheet: row 3 is b.md ---
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-19 20:42:54"
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Range("P1").EntireColumn.ColumnWidth = 39.1666666666667
